$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows at 727-732, pushing the existing rows 727-740
# down to 733-746 (this also grows the sheet dimension from T740 to T746).
$ws.Rows("727:732").Insert()

# Common (unchanged) column values shared by every row in this block.
$mercadoId = 3
$mercado = 'Femacal de La Calera'
$region = 'Coquimbo'
$codreg = 5
$tipo = 'Fruta'
$productoId = 100103
$producto = 'Frutos de hueso (carozo)'
$categoriaId = 100103006
$categoria = 'Nectarín'

# New data for rows 727-732 (June Pearl / Venus varieties, fecha 2022-02-03).
$rows = @(
    @{ Row=727; Fecha=44595; Variedad='June Pearl'; Calidad='Especial'; Volumen=85; PMin=16000; PMax=16000; PProm=16000; Unidad='$/caja 15 kilos empedrada'; Origen="Región de O'Higgins"; PKg=1067; KgUnidad=15 },
    @{ Row=728; Fecha=44595; Variedad='June Pearl'; Calidad='Primera';  Volumen=87; PMin=14000; PMax=14000; PProm=14000; Unidad='$/caja 15 kilos empedrada'; Origen="Región de O'Higgins"; PKg=933;  KgUnidad=15 },
    @{ Row=729; Fecha=44595; Variedad='June Pearl'; Calidad='Segunda';  Volumen=80; PMin=12000; PMax=12000; PProm=12000; Unidad='$/caja 15 kilos empedrada'; Origen="Región de O'Higgins"; PKg=800;  KgUnidad=15 },
    @{ Row=730; Fecha=44595; Variedad='Venus';      Calidad='Especial'; Volumen=90; PMin=16000; PMax=16000; PProm=16000; Unidad='$/caja 15 kilos empedrada'; Origen="Región de O'Higgins"; PKg=1067; KgUnidad=15 },
    @{ Row=731; Fecha=44595; Variedad='Venus';      Calidad='Primera';  Volumen=97; PMin=14000; PMax=14000; PProm=14000; Unidad='$/caja 15 kilos empedrada'; Origen="Región de O'Higgins"; PKg=933;  KgUnidad=15 },
    @{ Row=732; Fecha=44595; Variedad='Venus';      Calidad='Segunda';  Volumen=90; PMin=12000; PMax=12000; PProm=12000; Unidad='$/caja 15 kilos empedrada'; Origen="Región de O'Higgins"; PKg=800;  KgUnidad=15 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
